$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("T148").NumberFormat = "0.0000"
Write-Output "OK"
